# Updated cryptos list on Tue Jun 27 03:52:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain text values (e.g. "30.386.31", "1.001"),
# not numbers. Excel auto-converts numeric-looking strings assigned via
# .Value, so force the cell format to Text first to keep them as strings,
# matching the original workbook's inline-string cell type.

# --- Rows whose Price / Volume(1h) values change, ranking order unchanged ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.386.31"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.872.63"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.54"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2817"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06536"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.878.87"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07445"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.45"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.079"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.96"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6559"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.398.76"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.30"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007627"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.126.48"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.292"
$ws.Range("E21").Value = "  -0.45%  "

# --- Rows 22 & 23 swap rank order (BitcoinCash <-> BinanceUSD) ---
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.29"
$ws.Range("E23").Value = "  +12.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.188"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.242"
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.03"
$ws.Range("E26").Value = "  +4.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.54"
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.977"
$ws.Range("E28").Value = "  +3.26%  "
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09400"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.299"
$ws.Range("E31").Value = "  +0.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05036"
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.212"
$ws.Range("E34").Value = "  +10.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7552"
$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01836"
$ws.Range("E37").Value = "  +2.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.621"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.080"
$ws.Range("E39").Value = "  +1.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9064"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.941"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.74"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4285"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.443"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.55"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1301"
$ws.Range("E47").Value = "  -1.35%  "

# --- Rows 48 & 49 swap rank order (NEARProtocol <-> EnergySwap) ---
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.983"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.477"
$ws.Range("E49").Value = "  +7.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.16"
$ws.Range("E50").Value = "  +0.53%  "

# --- Row 51: Decentraland replaced by Cronos ---
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05653"
$ws.Range("E51").Value = "  -1.68%  "
